$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string "ECs" is inserted before "FAPs" in the shared string table.
# The data table below reflects the full, updated 9-row x 20-column dataset
# (sending/target clusters ECs, FAPs, sCs; ligand/receptor symbol Lrfn3).
$rows = @(
  @("ECs","Lrfn3","Lrfn3","ECs",2,0.6666666666666666,0.3670173333333334,1.101052,0.1110238110123159,0.1110238110123159,2,0.6666666666666666,0.3670173333333334,1.101052,0.1110238110123159,0.1110238110123159,0.1347017229671111,1.212315506704,0.01232628661169844,0.01232628661169844),
  @("ECs","Lrfn3","Lrfn3","FAPs",2,0.6666666666666666,0.3670173333333334,1.101052,0.1110238110123159,0.1110238110123159,3,1,2.43075,7.292249999999999,0.7353089462210328,0.7353089462210328,0.8921273829999999,8.029146447,0.08163680148090911,0.08163680148090913),
  @("ECs","Lrfn3","Lrfn3","sCs",2,0.6666666666666666,0.3670173333333334,1.101052,0.1110238110123159,0.1110238110123159,3,1,0.5079859999999999,1.523958,0.1536672427666513,0.1536672427666513,0.1864396670906667,1.677957003816,0.01706072291970836,0.01706072291970836),
  @("FAPs","Lrfn3","Lrfn3","ECs",3,1,2.43075,7.292249999999999,0.7353089462210328,0.7353089462210328,2,0.6666666666666666,0.3670173333333334,1.101052,0.1110238110123159,0.1110238110123159,0.8921273829999999,8.029146447,0.08163680148090911,0.08163680148090913),
  @("FAPs","Lrfn3","Lrfn3","FAPs",3,1,2.43075,7.292249999999999,0.7353089462210328,0.7353089462210328,3,1,2.43075,7.292249999999999,0.7353089462210328,0.7353089462210328,5.908545562499999,53.17691006249999,0.5406792463926857,0.5406792463926857),
  @("FAPs","Lrfn3","Lrfn3","sCs",3,1,2.43075,7.292249999999999,0.7353089462210328,0.7353089462210328,3,1,0.5079859999999999,1.523958,0.1536672427666513,0.1536672427666513,1.2347869695,11.1130827255,0.112992898347438,0.112992898347438),
  @("sCs","Lrfn3","Lrfn3","ECs",3,1,0.5079859999999999,1.523958,0.1536672427666513,0.1536672427666513,2,0.6666666666666666,0.3670173333333334,1.101052,0.1110238110123159,0.1110238110123159,0.1864396670906667,1.677957003816,0.01706072291970836,0.01706072291970836),
  @("sCs","Lrfn3","Lrfn3","FAPs",3,1,0.5079859999999999,1.523958,0.1536672427666513,0.1536672427666513,3,1,2.43075,7.292249999999999,0.7353089462210328,0.7353089462210328,1.2347869695,11.1130827255,0.112992898347438,0.112992898347438),
  @("sCs","Lrfn3","Lrfn3","sCs",3,1,0.5079859999999999,1.523958,0.1536672427666513,0.1536672427666513,3,1,0.5079859999999999,1.523958,0.1536672427666513,0.1536672427666513,0.2580497761959999,2.322447985764,0.02361362149950493,0.02361362149950494)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowdata = $rows[$i]
    for ($j = 0; $j -lt $rowdata.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $rowdata[$j]
    }
}
